$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.311.12"
$ws.Range("E2").Value = "  +2.48%  "
$ws.Range("D3").Value = "2.973.51"
$ws.Range("E3").Value = "  +1.33%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.66"
$ws.Range("E5").Value = "  +1.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.26"
$ws.Range("E6").Value = "  +4.00%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.520"
$ws.Range("E8").Value = "  +1.01%  "
$ws.Range("D9").Value = "2.965.94"
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.133"
$ws.Range("E10").Value = "  +4.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.34"
$ws.Range("E11").Value = "  +12.09%  "
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000229"
$ws.Range("E13").Value = "  +3.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.60"
$ws.Range("E14").Value = "  +2.62%  "
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "3.457.30"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.07"
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("D18").Value = "2.969.41"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "59.309.89"
$ws.Range("E19").Value = "  +2.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "436.14"
$ws.Range("E20").Value = "  +5.02%  "
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.718"
$ws.Range("E22").Value = "  +3.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.33"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.01"
$ws.Range("E24").Value = "  +0.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.69"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.23"
$ws.Range("E27").Value = "  +11.01%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.54"
$ws.Range("E29").Value = "  +2.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.72"
$ws.Range("E30").Value = "  +4.20%  "
$ws.Range("E31").Value = "  +8.05%  "
$ws.Range("E32").Value = "  +5.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.64"
$ws.Range("E33").Value = "  +1.02%  "
$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D34").Value = "0.0₃0765"
$ws.Range("E34").Value = "  +9.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.88"
$ws.Range("E35").Value = "  +3.89%  "
$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.983"
$ws.Range("E36").Value = "  +4.31%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.65"
$ws.Range("E39").Value = "  -1.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.76"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "399.17"
$ws.Range("E41").Value = "  +5.89%  "
$ws.Range("E42").Value = "  +1.75%  "
$ws.Range("D43").Value = "2.744.10"
$ws.Range("E43").Value = "  +1.79%  "
$ws.Range("E44").Value = "  -3.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.250"
$ws.Range("E45").Value = "  +6.42%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "34.57"
$ws.Range("E47").Value = "  +18.90%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.40"
$ws.Range("E48").Value = "  -1.14%  "
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("E50").Value = "  +2.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.28"
$ws.Range("E51").Value = "  +2.05%  "
